$wb = $excel.ActiveWorkbook

# --- "Games" sheet: append the completed game that was previously the
#     next upcoming fixture (BRK on 45306) as row 41 ---
$games = $wb.Worksheets.Item("Games")

$games.Cells.Item(41, 1).Value = 40
$games.Cells.Item(41, 2).Value = 45306
$games.Cells.Item(41, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item(41, 3).Value = 3
$games.Cells.Item(41, 4).Value = 96
$games.Cells.Item(41, 5).Value = 92.40000000000001
$games.Cells.Item(41, 6).Value = 0.425
$games.Cells.Item(41, 7).Value = 10.9
$games.Cells.Item(41, 8).Value = 12.8
$games.Cells.Item(41, 9).Value = 0.253
$games.Cells.Item(41, 10).Value = 94.09999999999999
$games.Cells.Item(41, 11).Value = "BRK"
$games.Cells.Item(41, 12).Value = 95
$games.Cells.Item(41, 13).Value = 0.404
$games.Cells.Item(41, 14).Value = 7
$games.Cells.Item(41, 15).Value = 15.8
$games.Cells.Item(41, 16).Value = 0.152
$games.Cells.Item(41, 17).Value = 93.09999999999999
$games.Cells.Item(41, 18).Value = 0
$games.Cells.Item(41, 19).Value = 1

# --- "Next" sheet: that fixture is no longer "next" -- remove its row
#     (row 2) and let every following fixture shift up one row ---
$next = $wb.Worksheets.Item("Next")
$next.Rows.Item(2).Delete()

Write-Output "done"
